$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1198.94
$summary.Range("B4").Value = -1.06
$summary.Range("B5").Value = -0.61
$summary.Range("B6").Value = 35
$summary.Range("B7").Value = 15
$summary.Range("B9").Value = 42.86

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 98.94
$status.Range("D4").Value = 35
$status.Range("E4").Value = -1.06
$status.Range("F4").Value = -1.06
$status.Range("G4").Value = 42.86

# --- Append new trade row (#35) to "All Trades" and "MarketMaking" sheets ---
$newRow = @{
    A = 35
    B = "2026-02-17"
    C = "13:23:11"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.33
    G = 0.36
    H = "CLOSED"
    I = 9.0909
    J = 0.03
    K = 98.94
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.1
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowNum = 36

    $ws.Range("A$rowNum").Value = $newRow.A

    # Date-like text ("2026-02-17") is auto-detected as a real date by Excel,
    # so force the cell to text first, write it, then drop back to the
    # default style so no stray number-format survives on the cell.
    $ws.Range("B$rowNum").NumberFormat = "@"
    $ws.Range("B$rowNum").Value = $newRow.B
    $ws.Range("B$rowNum").Style = "Normal"

    $ws.Range("C$rowNum").Value = $newRow.C
    $ws.Range("D$rowNum").Value = $newRow.D
    $ws.Range("E$rowNum").Value = $newRow.E
    $ws.Range("F$rowNum").Value = $newRow.F
    $ws.Range("G$rowNum").Value = $newRow.G
    $ws.Range("H$rowNum").Value = $newRow.H
    $ws.Range("I$rowNum").Value = $newRow.I
    $ws.Range("J$rowNum").Value = $newRow.J
    $ws.Range("K$rowNum").Value = $newRow.K
    $ws.Range("L$rowNum").Value = $newRow.L
    $ws.Range("M$rowNum").Value = $newRow.M
    $ws.Range("N$rowNum").Value = $newRow.N
    $ws.Range("O$rowNum").Value = $newRow.O
    $ws.Range("P$rowNum").Value = $newRow.P
    $ws.Range("Q$rowNum").Value = $newRow.Q
}
